# Inventory Add functionality finished
#
# Appends 4 new WIP rows (92-95) to the tracking sheet. Every column in
# the existing data (dates, quantities, dollar amounts, NSNs, etc.) is
# stored as literal text rather than numbers/dates, so each new cell
# must defeat Excel's automatic type conversion the same way.
#
# Note: logic is written inline inside the loops (no helper function)
# because this COM host does not reliably apply property writes made on
# a COM Range object that was passed into a user-defined function and
# then invoked repeatedly from inside a loop.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 92; A = "4/18/2019"; B = "SPE5EM-19-V-4069"; C = "3";   D = "`$299.94 ";   E = "5330016570078"; F = "SEAL,PLAIN";          G = "KTSDI";         H = "180X210L08N"; I = "M33"; J = "2019 AUG 16" },
    @{ Row = 93; A = "4/18/2019"; B = "SPE7M5-19-V-8001"; C = "125"; D = "`$745.00";     E = "5920015723699"; F = "FUSE,CARTRIDGE";      G = "Rohde Schwarz"; H = "0099-6729-00"; I = "CP";  J = "2019 AUG 16" },
    @{ Row = 94; A = "4/18/2019"; B = "SPE7M2-19-V-1181"; C = "16";  D = "`$769.44 ";   E = "5355011119493"; F = "POINTER,DIAL";        G = "Cameron";       H = "0315-0005.B"; I = "M30"; J = "2019 SEP 05" },
    @{ Row = 95; A = "4/19/2019"; B = "SPE7MC-19-V-7094"; C = "3";   D = "`$4,609.86 "; E = "5915015956493"; F = "FILTER,RADIO FREQUE"; G = "Genisco";       H = "GF68200-50B"; I = "CP";  J = "2019 OCT 07" }
)

$cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")

foreach ($rowData in $newRows) {
    $r = $rowData.Row
    foreach ($col in $cols) {
        $rng = $ws.Range($col + $r)
        # Forcing Text format before the assignment stops Excel from
        # re-interpreting look-alike dates/numbers/currency as their
        # typed equivalents; resetting the style afterwards keeps the
        # cell from carrying a lingering "@" number-format style index,
        # matching the unstyled cells already used throughout the sheet.
        $rng.NumberFormat = "@"
        $rng.Value = $rowData[$col]
        $rng.Style = "Normal"
    }
}
